$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 159, shifting the existing rows 159..379 down to 160..380
$ws.Rows.Item(159).Insert()

# Populate the newly inserted row 159 with the new daily price record
$ws.Cells.Item(159, 1).Value = 3
$ws.Cells.Item(159, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(159, 3).Value = "Coquimbo"
$ws.Cells.Item(159, 4).Value = 44895
$ws.Cells.Item(159, 5).Value = 5
$ws.Cells.Item(159, 6).Value = 100112039
$ws.Cells.Item(159, 7).Value = "Ciboulette"
$ws.Cells.Item(159, 8).Value = "Sin especificar"
$ws.Cells.Item(159, 9).Value = "Primera"
$ws.Cells.Item(159, 10).Value = 120
$ws.Cells.Item(159, 11).Value = 1500
$ws.Cells.Item(159, 12).Value = 1500
$ws.Cells.Item(159, 13).Value = 1500
$ws.Cells.Item(159, 14).Value = "$/docena de atados"
$ws.Cells.Item(159, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(159, 16).Value = 500
$ws.Cells.Item(159, 17).Value = 3
$ws.Cells.Item(159, 18).Value = "Hortaliza"
